$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Diego Ramos Barros"
$ws.Range("C2").Value = 28

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "05/03/1995"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "Masculino"
$ws.Range("F2").Value = "Branco"
$ws.Range("G2").Value = "Médio completo"
$ws.Range("H2").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I2").Value = "(41) 1417-2308"
$ws.Range("J2").Value = "(41) 98806-7654"
$ws.Range("K2").Value = "099.959.210-65"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "694640277"
$ws.Range("L2").ClearFormats()
